# Update nombre_aides (column C) and montant_total (column E) values
# for the Fonds de solidarite 2022-06-15 data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 57;  C = 31596;   E = 162615330 },
    @{ Row = 91;  C = 18880;   E = 75348959 },
    @{ Row = 93;  C = 16946;   E = 50673087 },
    @{ Row = 96;  C = 29545;   E = 56437589 },
    @{ Row = 98;  C = 6300;    E = 19499128 },
    @{ Row = 115; C = 81810;   E = 436757089 },
    @{ Row = 121; C = 1306372; E = 2275411864 },
    @{ Row = 129; C = 633749;  E = 3434716128 },
    @{ Row = 132; C = 585997;  E = 3472457973 },
    @{ Row = 136; C = 26705;   E = 144407972 },
    @{ Row = 186; C = 236836;  E = 1189991300 },
    @{ Row = 189; C = 100472;  E = 556467136 },
    @{ Row = 196; C = 595497;  E = 984000022 },
    @{ Row = 215; C = 230261;  E = 408757751 },
    @{ Row = 237; C = 283324;  E = 1438438825 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
